$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Java Batch 2")

# Fill in the "Student" (D) column for existing rows 4-10.
# Order chosen to reproduce the shared-string table insertion order of the target file.
$ws.Range("D6").Value = "Andrew"
$ws.Range("D7").Value = "Yamini"
$ws.Range("D9").Value = "Archana"
$ws.Range("D10").Value = "Surya"
$ws.Range("D4").Value = "Subedha"
$ws.Range("D8").Value = "Sangeetha"
$ws.Range("D5").Value = "Meera"

# Add new row 11 for an additional project entry.
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Load Management for transport planningh"
$ws.Range("D11").Value = "Sivanesh"

$ws.Range("A11:D11").RowHeight = 30

# Update the view: scroll so row 10 area is visible and select D11, matching the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D11").Select()
